$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Portugal Primeira Liga")

# Row 139
$ws.Range("B139").Value2 = 6876581
$ws.Range("F139").Value2 = "Estrela"
$ws.Range("G139").Value2 = "Vizela"
$ws.Range("I139").Value2 = 1
$ws.Range("J139").Value2 = "D"
$ws.Range("K139").Value2 = 2
$ws.Range("L139").Value2 = 3.4
$ws.Range("M139").Value2 = 3.6
$ws.Range("N139").Value2 = 2.8
$ws.Range("O139").Value2 = 3.3
$ws.Range("P139").Value2 = 2.55
$ws.Range("Q139").Value2 = 0
$ws.Range("R139").Value2 = 2.025
$ws.Range("S139").Value2 = 1.825
$ws.Range("T139").Value2 = 2.25
$ws.Range("U139").Value2 = 1.975
$ws.Range("V139").Value2 = 1.875
$ws.Range("W139").Value2 = -1
$ws.Range("X139").Value2 = 2.3
$ws.Range("Z139").Value2 = 0
$ws.Range("AA139").Value2 = -0
$ws.Range("AB139").Value2 = -0.5
$ws.Range("AC139").Value2 = 0.4375

# Row 140
$ws.Range("B140").Value2 = 6876582
$ws.Range("F140").Value2 = "SC Farense"
$ws.Range("G140").Value2 = "Gil Vicente"
$ws.Range("I140").Value2 = 0
$ws.Range("J140").Value2 = "H"
$ws.Range("K140").Value2 = 1.8
$ws.Range("L140").Value2 = 3.6
$ws.Range("M140").Value2 = 4.2
$ws.Range("N140").Value2 = 2.2
$ws.Range("O140").Value2 = 3.25
$ws.Range("P140").Value2 = 3.25
$ws.Range("Q140").Value2 = -0.25
$ws.Range("R140").Value2 = 1.95
$ws.Range("S140").Value2 = 1.9
$ws.Range("T140").Value2 = 2.5
$ws.Range("U140").Value2 = 2.025
$ws.Range("V140").Value2 = 1.825
$ws.Range("W140").Value2 = 1.2
$ws.Range("X140").Value2 = -1
$ws.Range("Z140").Value2 = 0.95
$ws.Range("AA140").Value2 = -1
$ws.Range("AB140").Value2 = -1
$ws.Range("AC140").Value2 = 0.825

# Row 177
$ws.Range("B177").Value2 = 7758962
$ws.Range("F177").Value2 = "Estoril"
$ws.Range("G177").Value2 = "Estrela"
$ws.Range("H177").Value2 = 1
$ws.Range("K177").Value2 = 2.2
$ws.Range("L177").Value2 = 3.4
$ws.Range("M177").Value2 = 3.2
$ws.Range("N177").Value2 = 2.05
$ws.Range("O177").Value2 = 3.6
$ws.Range("P177").Value2 = 3.5
$ws.Range("Q177").Value2 = -0.5
$ws.Range("R177").Value2 = 2.06
$ws.Range("S177").Value2 = 1.84
$ws.Range("T177").Value2 = 2.5
$ws.Range("U177").Value2 = 1.975
$ws.Range("V177").Value2 = 1.875
$ws.Range("W177").Value2 = 1.05
$ws.Range("Z177").Value2 = 1.06
$ws.Range("AB177").Value2 = -1
$ws.Range("AC177").Value2 = 0.875

# Row 178
$ws.Range("B178").Value2 = 7758961
$ws.Range("F178").Value2 = "Benfica"
$ws.Range("G178").Value2 = "Gil Vicente"
$ws.Range("H178").Value2 = 3
$ws.Range("K178").Value2 = 1.181
$ws.Range("L178").Value2 = 7
$ws.Range("M178").Value2 = 14
$ws.Range("N178").Value2 = 1.142
$ws.Range("O178").Value2 = 9
$ws.Range("P178").Value2 = 15
$ws.Range("Q178").Value2 = -2
$ws.Range("R178").Value2 = 1.825
$ws.Range("S178").Value2 = 2.025
$ws.Range("T178").Value2 = 3.25
$ws.Range("U178").Value2 = 1.825
$ws.Range("V178").Value2 = 2.025
$ws.Range("W178").Value2 = 0.1419999999999999
$ws.Range("Z178").Value2 = 0.825
$ws.Range("AB178").Value2 = -0.5
$ws.Range("AC178").Value2 = 0.5125

# Row 195
$ws.Range("B195").Value2 = 6876630
$ws.Range("F195").Value2 = "Benfica"
$ws.Range("G195").Value2 = "Vizela"
$ws.Range("H195").Value2 = 6
$ws.Range("I195").Value2 = 1
$ws.Range("J195").Value2 = "H"
$ws.Range("K195").Value2 = 1.111
$ws.Range("L195").Value2 = 8.5
$ws.Range("M195").Value2 = 21
$ws.Range("N195").Value2 = 1.1
$ws.Range("O195").Value2 = 9.5
$ws.Range("P195").Value2 = 23
$ws.Range("Q195").Value2 = -2
$ws.Range("R195").Value2 = 1.89
$ws.Range("S195").Value2 = 2.01
$ws.Range("T195").Value2 = 3.25
$ws.Range("U195").Value2 = 1.825
$ws.Range("V195").Value2 = 2.025
$ws.Range("W195").Value2 = 0.1000000000000001
$ws.Range("Y195").Value2 = -1
$ws.Range("Z195").Value2 = 0.8899999999999999
$ws.Range("AA195").Value2 = -1
$ws.Range("AB195").Value2 = 0.825

# Row 196
$ws.Range("B196").Value2 = 6876633
$ws.Range("F196").Value2 = "Estoril"
$ws.Range("G196").Value2 = "Gil Vicente"
$ws.Range("H196").Value2 = 1
$ws.Range("I196").Value2 = 3
$ws.Range("J196").Value2 = "A"
$ws.Range("K196").Value2 = 2.2
$ws.Range("L196").Value2 = 3.4
$ws.Range("M196").Value2 = 3.2
$ws.Range("N196").Value2 = 2.1
$ws.Range("O196").Value2 = 3.5
$ws.Range("P196").Value2 = 3.5
$ws.Range("Q196").Value2 = -0.25
$ws.Range("R196").Value2 = 1.95
$ws.Range("S196").Value2 = 1.95
$ws.Range("T196").Value2 = 2.5
$ws.Range("U196").Value2 = 1.925
$ws.Range("V196").Value2 = 1.925
$ws.Range("W196").Value2 = -1
$ws.Range("Y196").Value2 = 2.5
$ws.Range("Z196").Value2 = -1
$ws.Range("AA196").Value2 = 0.95
$ws.Range("AB196").Value2 = 0.925

# Row 231
$ws.Range("B231").Value2 = 6875479
$ws.Range("F231").Value2 = "Moreirense"
$ws.Range("G231").Value2 = "Arouca"
$ws.Range("I231").Value2 = 0
$ws.Range("J231").Value2 = "H"
$ws.Range("K231").Value2 = 2.4
$ws.Range("L231").Value2 = 3.1
$ws.Range("M231").Value2 = 3.1
$ws.Range("N231").Value2 = 2.15
$ws.Range("P231").Value2 = 3.5
$ws.Range("Q231").Value2 = -0.25
$ws.Range("U231").Value2 = 1.975
$ws.Range("V231").Value2 = 1.875
$ws.Range("W231").Value2 = 1.15
$ws.Range("Y231").Value2 = -1
$ws.Range("Z231").Value2 = 0.8500000000000001
$ws.Range("AA231").Value2 = -1
$ws.Range("AB231").Value2 = -1
$ws.Range("AC231").Value2 = 0.875

# Row 232
$ws.Range("B232").Value2 = 6876663
$ws.Range("F232").Value2 = "Chaves"
$ws.Range("G232").Value2 = "Guimaraes"
$ws.Range("I232").Value2 = 2
$ws.Range("J232").Value2 = "A"
$ws.Range("K232").Value2 = 4.333
$ws.Range("L232").Value2 = 3.5
$ws.Range("M232").Value2 = 1.833
$ws.Range("N232").Value2 = 4
$ws.Range("P232").Value2 = 1.95
$ws.Range("Q232").Value2 = 0.5
$ws.Range("U232").Value2 = 1.875
$ws.Range("V232").Value2 = 1.975
$ws.Range("W232").Value2 = -1
$ws.Range("Y232").Value2 = 0.95
$ws.Range("Z232").Value2 = -1
$ws.Range("AA232").Value2 = 1
$ws.Range("AB232").Value2 = 0.875
$ws.Range("AC232").Value2 = -1

# Row 244
$ws.Range("N244").Value2 = 2.1
$ws.Range("O244").Value2 = 3.4
$ws.Range("Q244").Value2 = -0.25
$ws.Range("R244").Value2 = 1.82
$ws.Range("S244").Value2 = 2.08
$ws.Range("T244").Value2 = 2.25
$ws.Range("U244").Value2 = 1.825
$ws.Range("V244").Value2 = 2.025

# Row 245
$ws.Range("O245").Value2 = 3.3
$ws.Range("R245").Value2 = 2
$ws.Range("S245").Value2 = 1.9
$ws.Range("U245").Value2 = 2.05
$ws.Range("V245").Value2 = 1.8

# Row 247
$ws.Range("O247").Value2 = 4.75
$ws.Range("P247").Value2 = 7
$ws.Range("R247").Value2 = 1.98
$ws.Range("S247").Value2 = 1.92
$ws.Range("U247").Value2 = 1.825
$ws.Range("V247").Value2 = 2.025

# Row 248
$ws.Range("N248").Value2 = 1.95
$ws.Range("P248").Value2 = 3.75
$ws.Range("R248").Value2 = 1.99
$ws.Range("S248").Value2 = 1.91
$ws.Range("U248").Value2 = 1.875
$ws.Range("V248").Value2 = 1.975

# Row 249
$ws.Range("P249").Value2 = 3.4
$ws.Range("R249").Value2 = 1.91
$ws.Range("S249").Value2 = 1.99
$ws.Range("T249").Value2 = 2.25
$ws.Range("U249").Value2 = 1.85
$ws.Range("V249").Value2 = 2

# Row 250
$ws.Range("N250").Value2 = 2.05
$ws.Range("P250").Value2 = 3.75
$ws.Range("R250").Value2 = 2.06
$ws.Range("S250").Value2 = 1.84

# Row 251
$ws.Range("P251").Value2 = 7.5
$ws.Range("R251").Value2 = 1.86
$ws.Range("S251").Value2 = 2.04

# Row 252
$ws.Range("N252").Value2 = 2.55
$ws.Range("O252").Value2 = 3.1
$ws.Range("P252").Value2 = 3
$ws.Range("Q252").Value2 = 0
$ws.Range("R252").Value2 = 1.8
$ws.Range("S252").Value2 = 2.1
$ws.Range("U252").Value2 = 2.025
$ws.Range("V252").Value2 = 1.825
